$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1727.6578
$ws.Range("I137").Value = 1147.2632
$ws.Range("J137").Value = 2308.0527
$ws.Range("K137").Value = 3441.7896
$ws.Range("L137").Value = 6924.158100000001
$ws.Range("M137").Value = -891.7896000000001
$ws.Range("N137").Value = -12024.1581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3834.6738
$ws.Range("I61").Value = 4868.273
$ws.Range("K61").Value = 4868.273
$ws.Range("M61").Value = -4656.273

$ws.Range("H74").Value = 5169.6333
$ws.Range("I74").Value = 6667.8423
$ws.Range("J74").Value = 2581.818
$ws.Range("K74").Value = 6667.8423
$ws.Range("L74").Value = 2581.818
$ws.Range("M74").Value = -5793.8423
$ws.Range("N74").Value = -4329.818

$ws.Range("H77").Value = 5169.6333
$ws.Range("I77").Value = 6667.8423
$ws.Range("J77").Value = 2581.818
$ws.Range("K77").Value = 33339.2115
$ws.Range("L77").Value = 12909.09
$ws.Range("M77").Value = -28971.2115
$ws.Range("N77").Value = -21645.09

$ws.Range("H136").Value = 3834.6738
$ws.Range("I136").Value = 4868.273
$ws.Range("K136").Value = 14604.819
$ws.Range("M136").Value = -12054.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3773.2454
$ws.Range("I134").Value = 5233.037
$ws.Range("J134").Value = 2257.3076
$ws.Range("K134").Value = 15699.111
$ws.Range("L134").Value = 6771.9228
$ws.Range("M134").Value = -13164.111
$ws.Range("N134").Value = -11841.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2314.2588
$ws.Range("I31").Value = 1453.5
$ws.Range("J31").Value = 2614.8413
$ws.Range("K31").Value = 1453.5
$ws.Range("L31").Value = 2614.8413
$ws.Range("M31").Value = -1158.5
$ws.Range("N31").Value = -3204.8413

$ws.Range("H34").Value = 2314.2588
$ws.Range("I34").Value = 1453.5
$ws.Range("J34").Value = 2614.8413
$ws.Range("K34").Value = 1453.5
$ws.Range("L34").Value = 2614.8413
$ws.Range("M34").Value = -1251.5
$ws.Range("N34").Value = -3018.8413

$ws.Range("H58").Value = 1570.04
$ws.Range("I58").Value = 1210.6
$ws.Range("J58").Value = 2109.2
$ws.Range("K58").Value = 1210.6
$ws.Range("L58").Value = 2109.2
$ws.Range("M58").Value = -1007.6
$ws.Range("N58").Value = -2515.2

$ws.Range("H132").Value = 3819.077
$ws.Range("I132").Value = 3105.1428
$ws.Range("J132").Value = 4652
$ws.Range("K132").Value = 9315.428400000001
$ws.Range("L132").Value = 13956
$ws.Range("M132").Value = -6785.428400000001
$ws.Range("N132").Value = -19016

$ws.Range("H134").Value = 2897.2964
$ws.Range("I134").Value = 3138.625
$ws.Range("J134").Value = 966.6667
$ws.Range("K134").Value = 9415.875
$ws.Range("L134").Value = 2900.0001
$ws.Range("M134").Value = -6880.875
$ws.Range("N134").Value = -7970.0001

$ws.Range("H136").Value = 1570.04
$ws.Range("I136").Value = 1210.6
$ws.Range("J136").Value = 2109.2
$ws.Range("K136").Value = 3631.8
$ws.Range("L136").Value = 6327.599999999999
$ws.Range("M136").Value = -1081.8
$ws.Range("N136").Value = -11427.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 53001.668
$ws.Range("J48").Value = 53001.668
$ws.Range("L48").Value = 159005.004
$ws.Range("N48").Value = -159505.004

$ws.Range("H54").Value = 6750
$ws.Range("J54").Value = 6750
$ws.Range("L54").Value = 20250
$ws.Range("N54").Value = -21368

$ws.Range("H63").Value = 3444.4443
$ws.Range("J63").Value = 6000
$ws.Range("L63").Value = 18000
$ws.Range("N63").Value = -19498

$ws.Range("H66").Value = 3444.4443
$ws.Range("J66").Value = 6000
$ws.Range("L66").Value = 54000
$ws.Range("N66").Value = -61488

$ws.Range("H68").Value = 2771.3108
$ws.Range("J68").Value = 2055.4285
$ws.Range("L68").Value = 6166.2855
$ws.Range("N68").Value = -7788.2855

$ws.Range("H71").Value = 2771.3108
$ws.Range("J71").Value = 2055.4285
$ws.Range("L71").Value = 18498.8565
$ws.Range("N71").Value = -26610.8565

$ws.Range("H75").Value = 28577230
$ws.Range("I75").Value = 5013
$ws.Range("J75").Value = 35720284
$ws.Range("K75").Value = 15039
$ws.Range("L75").Value = 107160852
$ws.Range("M75").Value = -14041
$ws.Range("N75").Value = -107162848

$ws.Range("H78").Value = 28577230
$ws.Range("I78").Value = 5013
$ws.Range("J78").Value = 35720284
$ws.Range("K78").Value = 45117
$ws.Range("L78").Value = 321482556
$ws.Range("M78").Value = -40125
$ws.Range("N78").Value = -321492540

$ws.Range("H94").Value = 3626.9285
$ws.Range("I94").Value = 175
$ws.Range("J94").Value = 4202.25
$ws.Range("K94").Value = 525
$ws.Range("L94").Value = 12606.75
$ws.Range("M94").Value = 151
$ws.Range("N94").Value = -13958.75

$ws.Range("H96").Value = 6370.3335
$ws.Range("I96").Value = 2026
$ws.Range("J96").Value = 7239.2
$ws.Range("K96").Value = 6078
$ws.Range("L96").Value = 21717.6
$ws.Range("M96").Value = -4019
$ws.Range("N96").Value = -25835.6

$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -566

$ws.Range("H107").Value = 1144.4423
$ws.Range("I107").Value = 336.4
$ws.Range("K107").Value = 1009.2
$ws.Range("M107").Value = 910.8000000000001

$ws.Range("H109").Value = 1942.3077
$ws.Range("J109").Value = 3283.3333
$ws.Range("L109").Value = 9849.999899999999
$ws.Range("N109").Value = -11929.9999

$ws.Range("H110").Value = 3669
$ws.Range("I110").Value = 3003.5
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 9010.5
$ws.Range("L110").Value = 15000
$ws.Range("M110").Value = -4920.5
$ws.Range("N110").Value = -23180

$ws.Range("H119").Value = 218586
$ws.Range("I119").Value = 8000
$ws.Range("J119").Value = 288781.34
$ws.Range("K119").Value = 24000
$ws.Range("L119").Value = 866344.02
$ws.Range("M119").Value = -19162
$ws.Range("N119").Value = -876020.02

$ws.Range("H134").Value = 10687.229
$ws.Range("I134").Value = 12455.3
$ws.Range("J134").Value = 9980
$ws.Range("K134").Value = 37365.89999999999
$ws.Range("L134").Value = 29940
$ws.Range("M134").Value = -32295.89999999999
$ws.Range("N134").Value = -40080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4587.125
$ws.Range("I126").Value = 5282.4614
$ws.Range("J126").Value = 1574
$ws.Range("K126").Value = 15847.3842
$ws.Range("L126").Value = 4722
$ws.Range("M126").Value = -13377.3842
$ws.Range("N126").Value = -9662

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 125003430
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 142860640
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 142860640
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -142862138

$ws.Range("H71").Value = 125003430
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 142860640
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 714303200
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -714310688

$ws.Range("H132").Value = 17340032
$ws.Range("I132").Value = 43345776
$ws.Range("J132").Value = 2869.7334
$ws.Range("K132").Value = 130037328
$ws.Range("L132").Value = 8609.200199999999
$ws.Range("M132").Value = -130034798
$ws.Range("N132").Value = -13669.2002

$ws.Range("H136").Value = 6016.1055
$ws.Range("I136").Value = 4705.7837
$ws.Range("J136").Value = 8440.200000000001
$ws.Range("K136").Value = 14117.3511
$ws.Range("L136").Value = 25320.6
$ws.Range("M136").Value = -11567.3511
$ws.Range("N136").Value = -30420.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1700
$ws.Range("I81").Value = 1700
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3400
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -2339

$ws.Range("H84").Value = 1700
$ws.Range("I84").Value = 1700
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 17000
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -11696

$ws.Range("H132").Value = 1953.6809
$ws.Range("I132").Value = 1431.3214
$ws.Range("J132").Value = 2723.4736
$ws.Range("K132").Value = 4293.9642
$ws.Range("L132").Value = 8170.4208
$ws.Range("M132").Value = -1763.9642
$ws.Range("N132").Value = -13230.4208

$ws.Range("H136").Value = 1294.8846
$ws.Range("I136").Value = 902.8889
$ws.Range("J136").Value = 2176.875
$ws.Range("K136").Value = 2708.6667
$ws.Range("L136").Value = 6530.625
$ws.Range("M136").Value = -158.6667000000002
$ws.Range("N136").Value = -11630.625
